# Applies the two changes captured by the commit:
#
#  1. Slide 5's table ("Google Shape;122;p17") is re-styled: its
#     <a:tableStyleId> moves from {630D055C-E541-4B05-9439-583696817DAE}
#     to {5200B4C3-5C0C-40DE-A88C-7D170456FDA9} (PowerPoint's "Medium
#     Style 2 - Accent 1" built-in table style GUID).
#
#  2. The deck's theme palette is switched from the custom "Integral /
#     Red Violet" scheme over to the stock Office theme palette.

$p = $ppt.ActivePresentation

# -- 1. Table style -----------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{5200B4C3-5C0C-40DE-A88C-7D170456FDA9}")

# -- 2. Theme colours -----------------------------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor $colors 1  "000000"   # dk1
Set-ThemeColor $colors 2  "FFFFFF"   # lt1
Set-ThemeColor $colors 3  "44546A"   # dk2
Set-ThemeColor $colors 4  "E7E6E6"   # lt2
Set-ThemeColor $colors 5  "5B9BD5"   # accent1
Set-ThemeColor $colors 6  "ED7D31"   # accent2
Set-ThemeColor $colors 7  "A5A5A5"   # accent3
Set-ThemeColor $colors 8  "FFC000"   # accent4
Set-ThemeColor $colors 9  "4472C4"   # accent5
Set-ThemeColor $colors 10 "70AD47"   # accent6
Set-ThemeColor $colors 11 "0563C1"   # hlink
Set-ThemeColor $colors 12 "954F72"   # folHlink
